$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values
$ws.Range("B2").Value = 5.5
$ws.Range("C2").Value = 10.5
$ws.Range("B3").Value = 5
$ws.Range("C5").Value = 25

# Update column widths to match new bestFit sizes (AutoFit based on content)
$ws.Columns.Item(1).ColumnWidth = 20.714285714285715
$ws.Columns.Item(2).ColumnWidth = 4.428571428571429
$ws.Columns.Item(3).ColumnWidth = 4.714285714285714

$wb.Save()
